## preflight-prep-guide.docx edit script
## 1) Shrink the height of the "Operating Manual" table row.
## 2) Word's grammar checker flagged a handful of phrases ("gramStart"/
##    "gramEnd"); reproduce that by splitting the affected runs and
##    inserting <w:proofErr> markers around the flagged sub-phrase,
##    leaving the visible text identical.

$d = $word.ActiveDocument
$apos = [char]0x2019
$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Insert-ProofErrSplit($SearchText, $InnerXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $SearchText"
    }

    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document ' + $wNs + '>' +
        '<w:body>' +
        $InnerXml +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($pkg)
}

## --- 1. Row height: "Immediately available:" / "Operating Manual" row ---
$t = $d.Tables.Item(1)
$t.Rows.Item(5).Height = 18

## --- 2. "otherwise unfit to perform properly the person's duties" ---
Insert-ProofErrSplit "otherwise unfit to perform properly the person${apos}s duties" (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:ind w:left="528" w:hanging="168"/></w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>otherwise</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> unfit to perform properly the person' + $apos + 's duties</w:t></w:r>' +
    '</w:p>'
)

## --- 3. "ensure that there is a sufficient amount of fuel or energy..." ---
## (paragraph has 3 runs: "e" | "nsure...flight" | " (901.28)" -- InsertXML
## replaces the whole paragraph it lands in, so the search + replacement
## must cover the complete paragraph text, including the untouched runs,
## or the sibling runs get silently dropped.)
Insert-ProofErrSplit "ensure that there is a sufficient amount of fuel or energy for safe completion of the flight (901.28)" (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:ind w:left="528" w:hanging="168"/></w:pPr>' +
    '<w:r><w:t>e</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">nsure that there is </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>a sufficient amount of</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> fuel or energy for safe completion of the flight</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (901.28)</w:t></w:r>' +
    '</w:p>'
)

## --- 4. "Take into account (901.27):" ---
Insert-ProofErrSplit "Take into account (901.27):" (
    '<w:p>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Take into account</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> (901.27):</w:t></w:r>' +
    '</w:p>'
)

## --- 5. "the proximity of aerodromes, airports and heliports" ---
Insert-ProofErrSplit "the proximity of aerodromes, airports and heliports" (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:ind w:left="528" w:hanging="168"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">the proximity of aerodromes, </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>airports</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and heliports</w:t></w:r>' +
    '</w:p>'
)

## --- 6. "there is no likelihood of collision with another aircraft, person or obstacle" ---
Insert-ProofErrSplit "there is no likelihood of collision with another aircraft, person or obstacle" (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:ind w:left="528" w:hanging="168"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">there is no likelihood of collision with another aircraft, </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>person</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> or obstacle</w:t></w:r>' +
    '</w:p>'
)

## --- 7. "the site set aside for take-off, launch, landing or recovery, as the case may be, is suitable for the intended" ---
## (paragraph also has a trailing Helvetica-styled space run + "operation" --
## include them unchanged for the same reason as item 3 above.)
Insert-ProofErrSplit "the site set aside for take-off, launch, landing or recovery, as the case may be, is suitable for the intended operation" (
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:ind w:left="528" w:hanging="168"/></w:pPr>' +
    '<w:r><w:t>the site set aside for take-off, launch, landing or recovery</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>, as the case may be, is</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> suitable for the intended</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>operation</w:t></w:r>' +
    '</w:p>'
)

Write-Output "edit complete"
